# Add the weekly figures for row 19 (date 2025-05-02, serial 45779) across
# the four tracking sheets, then leave each sheet's selection where the
# user last left it.

$wb = $excel.ActiveWorkbook

# --- Produzione ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Produzione")
$ws.Range("A19").Value = 45779
$ws.Range("B19").Value = 26459.32
$ws.Range("B19").NumberFormat = "0.00"
[void]$ws.Range("E18").Select()

# --- Entrate --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entrate")
$ws.Range("A19").Value = 45779
$ws.Range("B19").Value = 737.9

# --- Uscite -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Uscite")
$ws.Range("A19").Value = 45779
$ws.Range("B19").Value = 1456.19
$ws.Range("B19").NumberFormat = "0.00"
[void]$ws.Range("E17").Select()

# --- Saldo -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Saldo")
$ws.Range("A19").Value = 45779
$ws.Range("B19").Value = 56145.59
[void]$ws.Range("B22").Select()
